# Update Vcam1-Itga4 LR-pair sheet with newly computed TPM-based values.
# Columns (1-indexed): A=1 Sending cluster, D=4 Target cluster,
#   G=7,H=8,I=9,J=10 (Ligand avg/total expr + derived specificity, keyed by Sending cluster)
#   M=13,N=14,O=15,P=16 (Receptor avg/total expr + derived specificity, keyed by Target cluster)
#   Q=17,R=18,S=19,T=20 (Edge weights/specificity = products of the above)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand-side values (G,H,I,J) keyed by Sending cluster.
# $null entries mean that value is unchanged by this edit (only specificity columns moved).
$ghij = @{
    "ECs"               = @(22.95531766666667, 68.865953, 0.1720020945576478, 0.1720020945576478)
    "FAPs"              = @($null, $null, 0.4661646602805707, 0.4661646602805707)
    "Inflammatory-Mac"  = @(7.783044333333334, 23.349133, 0.05831763893698088, 0.05831763893698089)
    "MuSCs"             = @(30.44016466666666, 91.320494, 0.2280853681650076, 0.2280853681650076)
    "Resolving-Mac"     = @(10.06688366666667, 30.200651, 0.07543023805979308, 0.07543023805979308)
}

# New receptor-side values (M,N,O,P) keyed by Target cluster.
$mnop = @{
    "ECs"               = @(0.3331066666666667, 0.99932, 0.002125805913843485, 0.002125805913843485)
    "FAPs"              = @($null, $null, 0.0008775937418887864, 0.0008775937418887864)
    "Inflammatory-Mac"  = @(91.40156066666667, 274.204682, 0.5833025803538128, 0.5833025803538128)
    "MuSCs"             = @(0.5759770000000001, 1.727931, 0.00367574544541637, 0.00367574544541637)
    "Resolving-Mac"     = @(64.24849033333334, 192.745471, 0.4100182745450386, 0.4100182745450385)
}

for ($row = 2; $row -le 26; $row++) {
    $sending = $ws.Cells.Item($row, 1).Value2
    $target  = $ws.Cells.Item($row, 4).Value2

    $gRow = $ghij[$sending]
    $mRow = $mnop[$target]

    $g = $gRow[0]
    $h = $gRow[1]
    $i = $gRow[2]
    $j = $gRow[3]
    if ($null -eq $g) { $g = $ws.Cells.Item($row, 7).Value2 }
    if ($null -eq $h) { $h = $ws.Cells.Item($row, 8).Value2 }

    $m = $mRow[0]
    $n = $mRow[1]
    $o = $mRow[2]
    $p = $mRow[3]
    if ($null -eq $m) { $m = $ws.Cells.Item($row, 13).Value2 }
    if ($null -eq $n) { $n = $ws.Cells.Item($row, 14).Value2 }

    $ws.Cells.Item($row, 7).Value  = $g
    $ws.Cells.Item($row, 8).Value  = $h
    $ws.Cells.Item($row, 9).Value  = $i
    $ws.Cells.Item($row, 10).Value = $j

    $ws.Cells.Item($row, 13).Value = $m
    $ws.Cells.Item($row, 14).Value = $n
    $ws.Cells.Item($row, 15).Value = $o
    $ws.Cells.Item($row, 16).Value = $p

    $ws.Cells.Item($row, 17).Value = $g * $m
    $ws.Cells.Item($row, 18).Value = $h * $n
    $ws.Cells.Item($row, 19).Value = $i * $o
    $ws.Cells.Item($row, 20).Value = $j * $p
}
